$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    appears on Overview!E2:F2, E3:F3 and on the zh-cn / de-de status column C
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn handback info: Latest Target File (I), Latest Handback File (J),
#    Latest Handback DateTime (K)
# ---------------------------------------------------------------------------
$zhXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhDate = "2016-08-17 22:35:45"
$aUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8a1b7a54e682ba751164359b31bf1281f8d08ffd/e2e/a.md"
$bUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8a1b7a54e682ba751164359b31bf1281f8d08ffd/e2e/b.md"

$wsZh.Range("I2").Value = "a.md"
$wsZh.Range("J2").Value = $zhXlf
$wsZh.Range("K2").Value = $zhDate

$wsZh.Range("I3").Value = "a.md"
$wsZh.Range("J3").Value = $zhXlf
$wsZh.Range("K3").Value = $zhDate

# Rebuild hyperlinks in display order A2, I2, A3, I3 so relationship ids
# come out rId2..rId5 in that order (matches a regenerated report).
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $aUrl, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $aUrl, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $bUrl, "", "", "b.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $aUrl, "", "", "a.md")

# ---------------------------------------------------------------------------
# 3. de-de handback info
# ---------------------------------------------------------------------------
$deXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deDate = "2016-08-17 22:35:52"

$wsDe.Range("I2").Value = "a.md"
$wsDe.Range("J2").Value = $deXlf
$wsDe.Range("K2").Value = $deDate

$wsDe.Range("I3").Value = "a.md"
$wsDe.Range("J3").Value = $deXlf
$wsDe.Range("K3").Value = $deDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $aUrl, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $aUrl, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $bUrl, "", "", "b.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $aUrl, "", "", "a.md")

# ---------------------------------------------------------------------------
# 4. Column width adjustments to fit the longer text now present
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(10).ColumnWidth = 40

$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(10).ColumnWidth = 40
